# PTW: AXI4 Redesigned, Test added
#
# Sheet2 gets a new test-vector row (row 9 is corrected, rows 10-14 are new)
# describing IRQ/abort handling during the page-table-walk fetch, plus two
# columns (E "f2d" and I "Comment") are widened to fit the new text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet2 is already the active sheet in this workbook

# --- fix existing row 9 -----------------------------------------------
$ws.Range("C9").Value2 = "NONE"
$ws.Range("G9").Value2 = 1

# --- row 10 -------------------------------------------------------------
$ws.Range("A10").Value2 = "BRANCH_TAKEN"
$ws.Range("B10").Value2 = 1
$ws.Range("C10").Value2 = "NONE"
$ws.Range("D10").Value2 = "IDLE"
$ws.Range("E10").Value2 = "INTERRUPT_PENDING"
$ws.Range("F10").Value2 = "ready"
$ws.Range("G10").Value2 = 1
$ws.Range("H10").Value2 = 0

# --- row 11 -------------------------------------------------------------
$ws.Range("A11").Value2 = "BRANCH_TAKEN"
$ws.Range("B11").Value2 = 1
$ws.Range("C11").Value2 = "NONE"
$ws.Range("D11").Value2 = "IDLE"
$ws.Range("E11").Value2 = "NONE"
$ws.Range("F11").Value2 = "abort"
$ws.Range("G11").Value2 = 1
$ws.Range("I11").Value2 = "abort because IRQ handling"

# --- row 12 -------------------------------------------------------------
$ws.Range("A12").Value2 = "BRANCH_TAKEN"
$ws.Range("B12").Value2 = 1
$ws.Range("C12").Value2 = "MTVEC"
$ws.Range("D12").Value2 = "IDLE"
$ws.Range("E12").Value2 = "NONE"
$ws.Range("F12").Value2 = "branch taken"
$ws.Range("G12").Value2 = 0

# --- row 13 -------------------------------------------------------------
$ws.Range("A13").Value2 = "MTVEC"
$ws.Range("B13").Value2 = 0
$ws.Range("C13").Value2 = "NONE"
$ws.Range("D13").Value2 = "DONE"
$ws.Range("E13").Value2 = "INSTR"
$ws.Range("F13").Value2 = "not ready"
$ws.Range("G13").Value2 = 0

# --- row 14 -------------------------------------------------------------
$ws.Range("D14").Value2 = "IDLE"
$ws.Range("E14").Value2 = "INSTR"
$ws.Range("F14").Value2 = "ready"
$ws.Range("G14").Value2 = 0

# --- widen columns E and I to fit the new text (drop their old bestFit) -
$ws.Columns.Item(5).ColumnWidth = 21.035714285714285
$ws.Columns.Item(9).ColumnWidth = 23.660714285714285

# --- move the selection like the author left it -------------------------
$ws.Range("H12").Select() | Out-Null
